# Updated JPN model - 2025-08-07 13:22
# Rebuilds the ScenMap sheet with the new vervestacks scenario-mapping layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenMap")

# Start from a clean sheet - the old trastg~/V1G/V2G layout is fully replaced.
$ws.Cells.Clear()

# --- Row 1: scenario-group name prefixes, plus the "C" label over the group-by column ---
$ws.Range("A1").Value = "vstacks_t1~"
$ws.Range("B1").Value = "vstacks_t5~"
$ws.Range("C1").Value = "vstacks_w2~"
$ws.Range("H1").Value = "C"

# --- Row 2: the two sub-group names (driving the sg_ helper formulas in row 5) ---
$ws.Range("H2").Value = "ngfs"
$ws.Range("I2").Value = "timeslice"

# --- Row 4: section headers ---
$ws.Range("A4").Value = "~ScenMap"
$ws.Range("G4").Value = "~ScenG"

# --- Row 5: column headers + helper formulas ---
$ws.Range("A5").Value = "Oname"
$ws.Range("B5").Value = "Name"
$ws.Range("C5").Value = "Desc"
$ws.Range("D5").Value = "Ldesc"
$ws.Range("G5").Value = "Scen"
$ws.Range("H5").Formula = "=""sg_""&H2"
$ws.Range("I5").Formula = "=""sg_""&I2"

# --- Data blocks: 3 timeslice groups (3 days / 15 days / 2 weeks) x 7 scenarios ---
$scenarios = @("Delayed transition","Net Zero 2050","NDCs","Below 2deg","Current Policies","Low demand","Fragmented World")
$tsLabels  = @("3 days","15 days","2 weeks")
$tsSuffix  = @("_3d","_15d","_2w")
$rootRefs  = @('$A$1','$B$1','$C$1')

for ($block = 0; $block -lt 3; $block++) {
    $startRow = 6 + $block * 7
    for ($i = 0; $i -lt 7; $i++) {
        $row = $startRow + $i
        $ws.Range("A$row").Formula = "=" + $rootRefs[$block] + "&TEXT(N$row,""0000"")"
        $ws.Range("B$row").Formula = "=G$row"
        $ws.Range("G$row").Formula = "=H$row&P$row"
        if ($block -eq 0) {
            $ws.Range("H$row").Value = $scenarios[$i]
            $ws.Range("N$row").Value = $i + 1
        } else {
            $refRow = $row - 7
            $ws.Range("H$row").Formula = "=H$refRow"
            $ws.Range("N$row").Formula = "=N$refRow"
        }
        $ws.Range("I$row").Value = $tsLabels[$block]
        $ws.Range("P$row").Value = $tsSuffix[$block]
    }
}

# --- Sheet view: selection on ScenMap lands on E13 after the rebuild ---
$ws.Activate()
$ws.Range("E13").Select()

# --- Make TS_Defs the active tab (was "process map" before) ---
$tsDefs = $wb.Worksheets.Item("TS_Defs")
$tsDefs.Activate()
